$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.872.76"
$ws.Range("E2").Value = "  +2.29%  "
$ws.Range("D3").Value = "3.055.91"
$ws.Range("E3").Value = "  +2.29%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.16%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "526.16"
$ws.Range("E5").Value = "  +5.90%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "142.61"
$ws.Range("E6").Value = "  +5.49%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("E8").Value = "  +5.44%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "7.70"
$ws.Range("E9").Value = "  +6.84%  "
$ws.Range("E10").Value = "  +7.93%  "
$ws.Range("E11").Value = "  +6.21%  "
$ws.Range("E12").Value = "  +2.32%  "
$ws.Range("D13").Value = "3.578.02"
$ws.Range("E13").Value = "  +2.26%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.06"
$ws.Range("E14").Value = "  +7.96%  "
$ws.Range("E15").Value = "  +16.51%  "
$ws.Range("B16").Value = "WrappedBTC"
$ws.Range("C16").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D16").Value = "57.793.17"
$ws.Range("E16").Value = "  +2.10%  "
$ws.Range("B17").Value = "Polkadot"
$ws.Range("C17").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "6.28"
$ws.Range("E17").Value = "  +8.14%  "
$ws.Range("D18").Value = "3.051.21"
$ws.Range("E18").Value = "  +2.14%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.11"
$ws.Range("E19").Value = "  +5.94%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.16"
$ws.Range("E20").Value = "  +5.32%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "338.27"
$ws.Range("E21").Value = "  +4.18%  "
$ws.Range("E22").Value = "  +0.00%  "
$ws.Range("E23").Value = "  +7.44%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "65.05"
$ws.Range("E24").Value = "  +5.84%  "
$ws.Range("E25").Value = "  +6.48%  "
$ws.Range("D26").Value = "0.0₃0976"
$ws.Range("E26").Value = "  +8.45%  "
$ws.Range("E27").Value = "  +0.01%  "
$ws.Range("E28").Value = "  +6.35%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.38"
$ws.Range("E29").Value = "  +9.61%  "
$ws.Range("E30").Value = "  +6.63%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.24"
$ws.Range("E31").Value = "  +5.19%  "
$ws.Range("E32").Value = "  +5.00%  "
$ws.Range("B33").Value = "Monero"
$ws.Range("C33").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "156.43"
$ws.Range("E33").Value = "  +1.93%  "
$ws.Range("B34").Value = "NEARProtocol"
$ws.Range("C34").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.74"
$ws.Range("E34").Value = "  +6.01%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.00"
$ws.Range("E35").Value = "  +7.13%  "
$ws.Range("E36").Value = "  +3.30%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "26.24"
$ws.Range("E37").Value = "  +13.79%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0703"
$ws.Range("E38").Value = "  +3.90%  "
$ws.Range("D39").Value = "3.091.60"
$ws.Range("E39").Value = "  +2.29%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "37.80"
$ws.Range("E40").Value = "  +3.37%  "
$ws.Range("E41").Value = "  +8.86%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.999"
$ws.Range("E42").Value = "  -0.20%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.48"
$ws.Range("E43").Value = "  +5.36%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.663"
$ws.Range("E44").Value = "  +3.38%  "
$ws.Range("D45").Value = "2.325.18"
$ws.Range("E45").Value = "  +4.22%  "
$ws.Range("E46").Value = "  +3.59%  "
$ws.Range("B47").Value = "dogwifhat"
$ws.Range("C47").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.00"
$ws.Range("E47").Value = "  +2.65%  "
$ws.Range("B48").Value = "VeChain"
$ws.Range("C48").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0246"
$ws.Range("E48").Value = "  +3.27%  "
$ws.Range("B49").Value = "Cosmos"
$ws.Range("C49").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.05"
$ws.Range("E49").Value = "  +5.38%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "20.07"
$ws.Range("E50").Value = "  +5.50%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0902"
$ws.Range("E51").Value = "  +6.86%  "
